# Update database (slide the 5-year reporting window forward by one fiscal
# year: drop 1396/12, shift 1397..1400 left, add 1401/12) and refresh the
# read_price / change_price derived figures for the new "تاریخ انتشار" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (12 ماهه منتهی به ...) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (تاریخ انتشار) ---
$ws.Range("D9").Value = "1399-02-28 (9)"
$ws.Range("E9").Value = "1400-04-16 (8)"
$ws.Range("F9").Value = "1401-03-24 (9)"
$ws.Range("G9").Value = "1402-02-27 (7)"
# H9 is a bare "yyyy-mm-dd"-shaped string (no trailing "(n)" suffix), so a
# plain .Value assignment gets auto-recognized as a date literal and
# re-stamped with a date number format/new style. Force it to stay plain
# text: write it under a Text format, then restore the original row-9 cell
# formatting (borrowed from G9, which already carries that style) via a
# formats-only paste so the stored style index matches its neighbours.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-27"
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 806862
$ws.Range("E11").Value = 1153913
$ws.Range("F11").Value = 1008240
$ws.Range("G11").Value = 1632271
$ws.Range("H11").Value = 2999738

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) ---
$ws.Range("D12").Value = -654290
$ws.Range("E12").Value = -878387
$ws.Range("F12").Value = -747453
$ws.Range("G12").Value = -1139548
$ws.Range("H12").Value = -2387738

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 152572
$ws.Range("E13").Value = 275526
$ws.Range("F13").Value = 260787
$ws.Range("G13").Value = 492723
$ws.Range("H13").Value = 612000

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی ---
$ws.Range("D14").Value = -4953
$ws.Range("E14").Value = -19283
$ws.Range("F14").Value = -33021
$ws.Range("G14").Value = -60712
$ws.Range("H14").Value = -95353

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) ---
# D15 used to hold the literal text "-"; now it becomes a real numeric 0,
# matching the already-numeric cells to its right.
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 909
$ws.Range("E16").Value = 2495
$ws.Range("F16").Value = 4265
$ws.Range("G16").Value = 3929
$ws.Range("H16").Value = 34059

# --- Row 17: سود (زیان) عملیاتی ---
$ws.Range("D17").Value = 148528
$ws.Range("E17").Value = 258738
$ws.Range("F17").Value = 232031
$ws.Range("G17").Value = 435940
$ws.Range("H17").Value = 550706

# --- Row 18: هزینه های مالی ---
$ws.Range("D18").Value = -15773
$ws.Range("E18").Value = -15110
$ws.Range("F18").Value = -44626
$ws.Range("G18").Value = -168851
$ws.Range("H18").Value = -262925

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 720
$ws.Range("E19").Value = 706
$ws.Range("F19").Value = 1477
$ws.Range("G19").Value = 14983
$ws.Range("H19").Value = 29352

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 133475
$ws.Range("E20").Value = 244334
$ws.Range("F20").Value = 188882
$ws.Range("G20").Value = 282072
$ws.Range("H20").Value = 317133

# --- Row 21: مالیات ---
$ws.Range("D21").Value = -22119
$ws.Range("E21").Value = -33437
$ws.Range("F21").Value = -32796
$ws.Range("G21").Value = -57039
$ws.Range("H21").Value = -40052

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 111356
$ws.Range("E22").Value = 210897
$ws.Range("F22").Value = 156086
$ws.Range("G22").Value = 225033
$ws.Range("H22").Value = 277081

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (stays all zero) ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: سود (زیان) خالص ---
$ws.Range("D24").Value = 111356
$ws.Range("E24").Value = 210897
$ws.Range("F24").Value = 156086
$ws.Range("G24").Value = 225033
$ws.Range("H24").Value = 277081

# --- Row 25: سود هر سهم پس از کسر مالیات ---
$ws.Range("D25").Value = 1113
$ws.Range("E25").Value = 2108
$ws.Range("F25").Value = 201
$ws.Range("G25").Value = 290
$ws.Range("H25").Value = 358

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 100043
$ws.Range("E26").Value = 100043
$ws.Range("F26").Value = 775000
$ws.Range("G26").Value = 775000
$ws.Range("H26").Value = 775000

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ---
$ws.Range("D27").Value = 144
$ws.Range("E27").Value = 272
$ws.Range("F27").Value = 201
$ws.Range("G27").Value = 290
$ws.Range("H27").Value = 358
